# Update per-leve profit calculations across several profession sheets.
# (Sheets correspond to the game "Bahamut_Profits" report: ALC, ARM, BSM,
#  CRP, CUL, GSM, LTW, WVR -- market-price driven columns H:N are refreshed.)

$wb = $excel.ActiveWorkbook

function Set-Cells {
    param($SheetName, $Row, $Values)
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Cells.Item($Row, $col).Value = $Values[$col]
    }
}

# Column indices: H=8 I=9 J=10 K=11 L=12 M=13 N=14

# ---------------- ALC ----------------
Set-Cells "ALC" 8 @{8=225.55556; 10=375; 12=1125; 14=-1403}
Set-Cells "ALC" 12 @{8=86.5; 10=98.666664; 12=98.666664; 14=-438.666664}
Set-Cells "ALC" 62 @{8=77184.21000000001; 9=103927.9; 10=10325; 11=103927.9; 12=10325; 13=-103303.9; 14=-11573}
Set-Cells "ALC" 65 @{8=77184.21000000001; 9=103927.9; 10=10325; 11=519639.5; 12=51625; 13=-516519.5; 14=-57865}
Set-Cells "ALC" 98 @{8=4458.923; 9=3996; 10=6002; 11=3996; 12=6002; 13=-2498; 14=-8998}
Set-Cells "ALC" 112 @{8=1746.6666; 10=1861; 12=5583; 14=-7799}
Set-Cells "ALC" 122 @{8=4458.923; 9=3996; 10=6002; 11=11988; 12=18006; 13=-9538; 14=-22906}
Set-Cells "ALC" 123 @{8=25055.555; 10=25055.555; 12=25055.555; 14=-34855.555}

# ---------------- ARM ----------------
Set-Cells "ARM" 61 @{8=1126.2; 9=827.7273; 10=1947; 11=827.7273; 12=1947; 13=-615.7273; 14=-2371}
Set-Cells "ARM" 74 @{8=1360.6; 9=1473.3334; 10=909.6667; 11=1473.3334; 12=909.6667; 13=-599.3334; 14=-2657.6667}
Set-Cells "ARM" 77 @{8=1360.6; 9=1473.3334; 10=909.6667; 11=7366.666999999999; 12=4548.3335; 13=-2998.666999999999; 14=-13284.3335}
Set-Cells "ARM" 136 @{8=1126.2; 9=827.7273; 10=1947; 11=2483.1819; 12=5841; 13=66.81809999999996; 14=-10941}

# ---------------- CRP ----------------
Set-Cells "CRP" 31 @{8=2398.6938; 9=2496.5; 10=2060.818; 11=2496.5; 12=2060.818; 13=-2201.5; 14=-2650.818}
Set-Cells "CRP" 34 @{8=2398.6938; 9=2496.5; 10=2060.818; 11=2496.5; 12=2060.818; 13=-2294.5; 14=-2464.818}
Set-Cells "CRP" 99 @{8=3481.1428; 9=3190.2222; 11=3190.2222; 13=-1692.2222}
Set-Cells "CRP" 126 @{8=3481.1428; 9=3190.2222; 11=9570.6666; 13=-7100.6666}
Set-Cells "CRP" 132 @{8=3494.5264; 9=2350; 11=7050; 13=-4520}

# ---------------- CUL ----------------
Set-Cells "CUL" 94 @{8=3965.5386; 9=2140; 11=6420; 13=-5744}

# ---------------- GSM ----------------
Set-Cells "GSM" 113 @{8=1402.1538; 9=516.4286; 11=516.4286; 13=1653.5714}
Set-Cells "GSM" 122 @{8=3290475.5; 9=13157894; 10=1336; 11=39473682; 12=4008; 13=-39471232; 14=-8908}

# ---------------- LTW ----------------
Set-Cells "LTW" 40 @{8=1444102.9; 9=3367936.8; 10=1227.5; 11=3367936.8; 12=1227.5; 13=-3367800.8; 14=-1499.5}
Set-Cells "LTW" 46 @{8=2608; 9=2412.125; 10=2999.75; 11=2412.125; 12=2999.75; 13=-2224.125; 14=-3375.75}
Set-Cells "LTW" 119 @{8=45000; 10=45000; 12=45000; 14=-54676}
Set-Cells "LTW" 122 @{8=10014.385; 9=16385.715; 10=2581.1667; 11=49157.145; 12=7743.500100000001; 13=-46707.145; 14=-12643.5001}
Set-Cells "LTW" 136 @{8=3770.6; 9=1271.9286; 10=9600.833000000001; 11=3815.7858; 12=28802.499; 13=-1265.7858; 14=-33902.499}

# ---------------- WVR ----------------
Set-Cells "WVR" 122 @{8=1276.8572; 9=1207.4286; 10=1415.7142; 11=3622.2858; 12=4247.142599999999; 13=-1172.2858; 14=-9147.142599999999}
Set-Cells "WVR" 136 @{8=730.5714; 9=709.5; 10=758.6667; 11=2128.5; 12=2276.0001; 13=421.5; 14=-7376.0001}
